$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- long Cypher query strings (used verbatim, no PowerShell expansion) ----
$sCasesQuery = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
    WHERE a.arm_id IN ['Z1D']
OPTIONAL MATCH (f:file)-[*]->(c)
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@
$sFilesQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
WHERE a.arm_id IN ['Z1D']
WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@
$sStatQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
WHERE a.arm_id IN ['Z1D']
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

# ---- Insert a new first column for the "TabName" label; shifts old A:D -> B:E ----
$ws.Columns.Item(1).Insert()

# ---- Remember the (now-shifted) output-file names so row 3 can reuse them ----
$neo4jFile = $ws.Range("D2").Value2
$webFile   = $ws.Range("E2").Value2

# ---- Row 1 : header ----
$ws.Range("A1").Value = "TabName"

# ---- Row 2 : CasesTab ----
$ws.Range("A2").Value = "CasesTab"
$ws.Range("B2").Value = $sCasesQuery
$ws.Range("B2").WrapText = $true
$ws.Range("C2").Value = $sStatQuery
$ws.Range("C2").WrapText = $true

# ---- Row 3 : FilesTab (new row) ----
$ws.Range("A3").Value = "FilesTab"
$ws.Range("B3").Value = $sFilesQuery
$ws.Range("B3").WrapText = $true
$ws.Range("C3").Value = $sStatQuery
$ws.Range("C3").WrapText = $true
$ws.Range("D3").Value = $neo4jFile
$ws.Range("E3").Value = $webFile

# ---- Row heights sized for the wrapped query text ----
$ws.Rows.Item(2).RowHeight = 195
$ws.Rows.Item(3).RowHeight = 409.5

# ---- Auto-fit the new label column ----
$ws.Columns.Item(1).AutoFit()

# ---- Selection, matching the saved sheet view ----
$ws.Range("C2").Select() | Out-Null
